$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Productos")

# Update unit prices
$ws.Range("C7").Value = 175
$ws.Range("C8").Value = 175
$ws.Range("C9").Value = 190
$ws.Range("C10").Value = 190

# Update the active selection on the sheet
$ws.Activate()
$ws.Range("E22").Select()
